$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the applicant rows that were dropped from this Vacantes listing.
# Delete bottom-most rows first so earlier row numbers stay valid.
$ws.Rows(47).Delete()   # nro_inscripcion 151 - Diego Torres
$ws.Rows(46).Delete()   # nro_inscripcion 147 - Alberto Ramos
$ws.Rows(42).Delete()   # nro_inscripcion 123 - Fernando Peralta
$ws.Rows(20).Delete()   # nro_inscripcion 36  - Luisa Aguilera

# Insert a new applicant row ("Luisa Pérez") ahead of nro_inscripcion 92 (row 36 now).
$ws.Rows(36).Insert()
$ws.Range("A36").Value = 88
$ws.Range("B36").Value = "Luisa Pérez"
$ws.Range("C36").Value = 9
$ws.Range("D36").Value = 4

# Replace what is now the final row (previously nro_inscripcion 172 - Julieta Sánchez)
# with the new last entrant.
$ws.Range("A47").Value = 200
$ws.Range("B47").Value = "Cecilia Aguilera"
$ws.Range("C47").Value = 10
$ws.Range("D47").Value = 5
